$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data refresh (GitHub Actions scheduled update).
# A helper writes each cell as literal text: NumberFormat is switched to "@"
# before the assignment so numeric-looking strings (e.g. "1.013",
# "28.613.10") are stored verbatim rather than being parsed into numbers,
# then the style is reset to "Normal" so no stray formatting sticks around.
function Set-TextCell([string]$ref, [string]$value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '28.613.10'
Set-TextCell 'E2' '  -1.56%  '
Set-TextCell 'D3' '1.968.35'
Set-TextCell 'E4' '  +0.40%  '
Set-TextCell 'D5' '323.73'
Set-TextCell 'E5' '  -0.88%  '
Set-TextCell 'D6' '1.012'
Set-TextCell 'E6' '  +0.47%  '
Set-TextCell 'D7' '0.4820'
Set-TextCell 'E7' '  -3.23%  '
Set-TextCell 'D8' '0.4078'
Set-TextCell 'E8' '  -3.17%  '
Set-TextCell 'D9' '54.20'
Set-TextCell 'E9' '  -0.17%  '
Set-TextCell 'D10' '0.08534'
Set-TextCell 'E10' '  -5.50%  '
Set-TextCell 'D11' '1.065'
Set-TextCell 'E11' '  -2.99%  '
Set-TextCell 'D12' '22.52'
Set-TextCell 'E12' '  -2.20%  '
Set-TextCell 'D13' '1.996.84'
Set-TextCell 'E13' '  +3.20%  '
Set-TextCell 'D14' '7.644'
Set-TextCell 'E14' '  -2.89%  '
Set-TextCell 'D15' '6.213'
Set-TextCell 'E15' '  -3.51%  '
Set-TextCell 'D16' '1.013'
Set-TextCell 'E16' '  +0.47%  '
Set-TextCell 'D17' '91.35'
Set-TextCell 'E17' '  +0.10%  '
Set-TextCell 'D18' '0.00001077'
Set-TextCell 'E18' '  -1.78%  '
Set-TextCell 'D19' '0.06659'
Set-TextCell 'E19' '  -0.14%  '
Set-TextCell 'D20' '18.67'
Set-TextCell 'E20' '  -2.48%  '
Set-TextCell 'D21' '1.012'
Set-TextCell 'E21' '  +0.55%  '
Set-TextCell 'D22' '5.887'
Set-TextCell 'E22' '  -0.92%  '
Set-TextCell 'D23' '28.642.95'
Set-TextCell 'D24' '11.55'
Set-TextCell 'E24' '  -2.99%  '
Set-TextCell 'D25' '2.304'
Set-TextCell 'E25' '  +0.82%  '
Set-TextCell 'D26' '2.224.88'
Set-TextCell 'E26' '  +2.57%  '
Set-TextCell 'D27' '156.69'
Set-TextCell 'E27' '  +0.55%  '
Set-TextCell 'D28' '20.43'
Set-TextCell 'E28' '  -0.80%  '
Set-TextCell 'D29' '5.930'
Set-TextCell 'E29' '  -3.62%  '
Set-TextCell 'D30' '2.190'
Set-TextCell 'E30' '  -2.85%  '
Set-TextCell 'D31' '125.13'
Set-TextCell 'E31' '  -1.67%  '
Set-TextCell 'D32' '0.9951'
Set-TextCell 'E32' '  -4.33%  '
Set-TextCell 'D33' '0.09696'
Set-TextCell 'E33' '  -1.38%  '
Set-TextCell 'D34' '1.466'
Set-TextCell 'E34' '  -4.27%  '
Set-TextCell 'D35' '3.711'
Set-TextCell 'E35' '  +0.95%  '
Set-TextCell 'D36' '5.663'
Set-TextCell 'E36' '  -2.25%  '
Set-TextCell 'D37' '9.235'
Set-TextCell 'E37' '  +3.38%  '
Set-TextCell 'D38' '0.02346'
Set-TextCell 'E38' '  -2.97%  '
Set-TextCell 'D39' '0.06262'
Set-TextCell 'E39' '  -0.35%  '
Set-TextCell 'D40' '1.261'
Set-TextCell 'E40' '  -2.22%  '
Set-TextCell 'D41' '0.6261'
Set-TextCell 'E41' '  -2.82%  '
Set-TextCell 'D42' '11.28'
Set-TextCell 'E42' '  -1.59%  '
Set-TextCell 'E43' '  +0.51%  '
Set-TextCell 'E44' '  -2.91%  '
Set-TextCell 'D45' '1.355'
Set-TextCell 'E45' '  +6.45%  '
Set-TextCell 'B46' 'EnergySwap'
Set-TextCell 'C46' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D46' '13.15'
Set-TextCell 'E46' '  -1.55%  '
Set-TextCell 'B47' 'Decentraland'
Set-TextCell 'C47' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextCell 'D47' '0.5986'
Set-TextCell 'E47' '  -3.28%  '
Set-TextCell 'D48' '2.079'
Set-TextCell 'E48' '  -4.03%  '
Set-TextCell 'D49' '3.420'
Set-TextCell 'E49' '  -1.18%  '
Set-TextCell 'D50' '0.06848'
Set-TextCell 'E50' '  -0.54%  '
Set-TextCell 'B51' 'BabyDogeCoin'
Set-TextCell 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 'D51' '0.00000000315'
Set-TextCell 'E51' '  -4.11%  '
